$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variable_mapping")

$ws.Range("A1").Value = "PlatformName"
$ws.Range("B1").Value = "HeaderRow"
